$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# C1 gains the "text" style already used by B1, G1..M1 (numFmt "@")
$ws.Cells.Item(1,3).NumberFormat = "@"
# N1 ("file_number") keeps its text, nothing else to do there.

# --- Row 2 : Cariboo / 240263 ---
$ws.Cells.Item(2,1).Value = "Cariboo"
$ws.Cells.Item(2,2).ClearContents()
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "0240263"
$ws.Cells.Item(2,4).Value = 949613
$ws.Cells.Item(2,5).Value = 984515
$ws.Cells.Item(2,6).Value = "\\spatialfiles.bcgov\Work\srm\wml\Workarea\Authorizations\Land\Cariboo\Cariboo_Batch_GR_20204_724\240263"
# G2:M2 already hold "False"/"False"/"False"/"False"/"False"/"True"/"False" - leave as-is.
$ws.Cells.Item(2,14).ClearContents()

# --- Row 3 : Cariboo / 247611 ---
$ws.Cells.Item(3,1).Value = "Cariboo"
$ws.Cells.Item(3,2).ClearContents()
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "0247611"
$ws.Cells.Item(3,4).Value = 949614
$ws.Cells.Item(3,5).Value = 984516
$ws.Cells.Item(3,6).Value = "T:\job2\\spatialfiles.bcgov\Work\srm\wml\Workarea\Authorizations\Land\Cariboo\Cariboo_Batch_GR_20204_724\247611"
# G3:M3 already hold "False"/"False"/"False"/"False"/"False"/"True"/"False" - leave as-is.
$ws.Cells.Item(3,14).ClearContents()

# --- Row 4 (new) : Cariboo / 5407734 ---
$ws.Cells.Item(4,1).Value = "Cariboo"
$ws.Cells.Item(4,3).Value = 5407734
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = 949615
$ws.Cells.Item(4,5).Value = 984517
$ws.Cells.Item(4,6).Value = "\\spatialfiles.bcgov\Work\srm\wml\Workarea\Authorizations\Land\Cariboo\Cariboo_Batch_GR_20204_724\25407734"
$ws.Range("G2:M2").Copy()
$ws.Range("G4:M4").PasteSpecial(-4163)

# --- Row 5 (new) : Cariboo / 5407781 ---
$ws.Cells.Item(5,1).Value = "Cariboo"
$ws.Cells.Item(5,3).Value = 5407781
$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = 949616
$ws.Cells.Item(5,5).Value = 984518
$ws.Cells.Item(5,6).Value = "\\spatialfiles.bcgov\Work\srm\wml\Workarea\Authorizations\Land\Cariboo\Cariboo_Batch_GR_20204_724\5407781"
$ws.Range("G2:M2").Copy()
$ws.Range("G5:M5").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Column widths (closest values reachable through ColumnWidth, which the
#     engine quantises to 1/6-character steps after MDW conversion) ---
$ws.Columns.Item(2).ColumnWidth = 31.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(6).ColumnWidth = 116.66666666666667

# --- Selection / view: scroll back to default and select C12 ---
$ws.Range("C12").Select()
